# Regenerate s_vals data to filter save games.
# Updates the numeric TB/d2S/K/IP/sum columns (B:E, G) for rows 2-12.
# Column F (Win) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 0.6753301551942219; C = 0.3127903958511391;  D = 3.900430680208489;  E = 0.496779210170732;  G = 5.385330441424582 }
    3  = @{ B = 1.459612070389937;  C = 1.667794583268128;   D = 0.1575252929769615; E = 0.496779210170732;  G = 3.781711156805759 }
    4  = @{ B = 3.230985683306322;  C = 1.667794583268128;   D = 3.900430680208489;  E = 0.496779210170732;  G = 9.295990156953671 }
    5  = @{ B = 3.230985683306322;  C = 1.667794583268128;   D = 0.1575252929769615; E = 0.496779210170732;  G = 5.553084769722144 }
    6  = @{ B = 0.6753301551942219; C = 1.667794583268128;   D = 0.8054896365839992; E = 0.496779210170732;  G = 3.645393585217082 }
    7  = @{ B = 0.6753301551942219; C = 0.3127903958511391;  D = 0.8054896365839992; E = 0.496779210170732;  G = 2.290389397800092 }
    8  = @{ B = 3.230985683306322;  C = 1.667794583268128;   D = 0.1575252929769615; E = 0.496779210170732;  G = 5.553084769722144 }
    9  = @{ B = 0.127881588408715;  C = 1.667794583268128;   D = 0.8054896365839992; E = 8.660232485948974;  G = 11.26139829420982 }
    10 = @{ B = 0.3048080303191223; C = 0.3127903958511391;  D = 0.8054896365839992; E = 0.496779210170732;  G = 1.919867272924993 }
    11 = @{ B = 3.230985683306322;  C = 1.667794583268128;   D = 26.21740644021617;  E = 8.660232485948974;  G = 39.7764191927396 }
    12 = @{ B = 0.6753301551942219; C = 0.3127903958511391;  D = 26.21740644021617;  E = 0.496779210170732;  G = 27.70230620143226 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("G$r").Value = $row.G
}
